$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New header cells (row 1): H1, I1, J1
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "Page Footer Links"
$ws.Range("I1").Value = "Impressum / Cookie / Privacy Policy Page footer"
$ws.Range("J1").Value = "Impressum length / horizontal scroll "

# ---------------------------------------------------------------------------
# Rows that simply get "+" in H, I and J
# ---------------------------------------------------------------------------
$plusRows = @(2,3,6,9,10,11,12,13,14,15,16,18,20,38,40)
foreach ($r in $plusRows) {
    $ws.Cells.Item($r, 8).Value = "+"
    $ws.Cells.Item($r, 9).Value = "+"
    $ws.Cells.Item($r, 10).Value = "+"
}

# ---------------------------------------------------------------------------
# Row 17 is special: I17 gets its own text instead of "+"
# ---------------------------------------------------------------------------
$ws.Range("H17").Value = "+"
$ws.Range("I17").Value = "cookie nicht rest +"
$ws.Range("J17").Value = "+"

# ---------------------------------------------------------------------------
# Row 19 is special: G19 is cleared, H19 and I19 get "+", no J19
# ---------------------------------------------------------------------------
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = "+"
$ws.Range("I19").Value = "+"

# ---------------------------------------------------------------------------
# Column widths for the new columns G:J (best-fit-like explicit widths)
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 13.893229166666666
$ws.Columns.Item(8).ColumnWidth = 14.166666666666666
$ws.Columns.Item(9).ColumnWidth = 38.346354166666664
$ws.Columns.Item(10).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# Scroll/selection: drop the old topLeftCell / selection, select J19 instead
# ---------------------------------------------------------------------------
$ws.Range("J19").Select() | Out-Null
